$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 53950
$ws.Range("J57").Value = 53950
$ws.Range("L57").Value = 161850
$ws.Range("N57").Value = -162848
$ws.Range("H112").Value = 5535.2705
$ws.Range("J112").Value = 6160.606
$ws.Range("L112").Value = 18481.818
$ws.Range("N112").Value = -20697.818
$ws.Range("H132").Value = 2336.6445
$ws.Range("I132").Value = 2143.0232
$ws.Range("K132").Value = 6429.069600000001
$ws.Range("M132").Value = -3899.069600000001
$ws.Range("H138").Value = 2450.09
$ws.Range("I138").Value = 1514.6154
$ws.Range("J138").Value = 2589.8735
$ws.Range("K138").Value = 4543.8462
$ws.Range("L138").Value = 7769.620500000001
$ws.Range("M138").Value = 596.1538
$ws.Range("N138").Value = -18049.6205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 20188.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 20188.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20188.5
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -20814.5
$ws.Range("H74").Value = 3440.6667
$ws.Range("I74").Value = 2412.6155
$ws.Range("J74").Value = 4655.636
$ws.Range("K74").Value = 2412.6155
$ws.Range("L74").Value = 4655.636
$ws.Range("M74").Value = -1538.6155
$ws.Range("N74").Value = -6403.636
$ws.Range("H77").Value = 3440.6667
$ws.Range("I77").Value = 2412.6155
$ws.Range("J77").Value = 4655.636
$ws.Range("K77").Value = 12063.0775
$ws.Range("L77").Value = 23278.18
$ws.Range("M77").Value = -7695.077499999999
$ws.Range("N77").Value = -32014.18
$ws.Range("H102").Value = 3333.3333
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -1378
$ws.Range("N102").Value = -6744
$ws.Range("H132").Value = 6849.3335
$ws.Range("I132").Value = 6189.2
$ws.Range("K132").Value = 18567.6
$ws.Range("M132").Value = -16037.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 92541.82000000001
$ws.Range("I107").Value = 144237.14
$ws.Range("J107").Value = 2075
$ws.Range("K107").Value = 144237.14
$ws.Range("L107").Value = 2075
$ws.Range("M107").Value = -142317.14
$ws.Range("N107").Value = -5915
$ws.Range("H134").Value = 3630.2
$ws.Range("I134").Value = 3564.9412
$ws.Range("K134").Value = 10694.8236
$ws.Range("M134").Value = -8159.8236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5976.4478
$ws.Range("I31").Value = 1778.4348
$ws.Range("J31").Value = 8170.864
$ws.Range("K31").Value = 1778.4348
$ws.Range("L31").Value = 8170.864
$ws.Range("M31").Value = -1483.4348
$ws.Range("N31").Value = -8760.864
$ws.Range("H34").Value = 5976.4478
$ws.Range("I34").Value = 1778.4348
$ws.Range("J34").Value = 8170.864
$ws.Range("K34").Value = 1778.4348
$ws.Range("L34").Value = 8170.864
$ws.Range("M34").Value = -1576.4348
$ws.Range("N34").Value = -8574.864
$ws.Range("H52").Value = 60780
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 60780
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 60780
$ws.Range("M52").Value = $null
$ws.Range("N52").Value = -61368
$ws.Range("H122").Value = 2019.3334
$ws.Range("I122").Value = 1800.75
$ws.Range("J122").Value = 2153.8462
$ws.Range("K122").Value = 5402.25
$ws.Range("L122").Value = 6461.5386
$ws.Range("M122").Value = -2952.25
$ws.Range("N122").Value = -11361.5386
$ws.Range("H139").Value = 39565.285
$ws.Range("J139").Value = 39565.285
$ws.Range("L139").Value = 39565.285
$ws.Range("N139").Value = -49845.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3433.3333
$ws.Range("I56").Value = 3433.3333
$ws.Range("K56").Value = 3433.3333
$ws.Range("M56").Value = -2903.3333
$ws.Range("H107").Value = 2550.4666
$ws.Range("I107").Value = 397.16666
$ws.Range("J107").Value = 3986
$ws.Range("K107").Value = 1191.49998
$ws.Range("L107").Value = 11958
$ws.Range("M107").Value = 728.5000199999999
$ws.Range("N107").Value = -15798
$ws.Range("H120").Value = 12500
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 12500
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 37500
$ws.Range("M120").Value = $null
$ws.Range("N120").Value = -47176
$ws.Range("H126").Value = 4025.9524
$ws.Range("I126").Value = 530
$ws.Range("J126").Value = 4200.75
$ws.Range("K126").Value = 1590
$ws.Range("L126").Value = 12602.25
$ws.Range("M126").Value = 3350
$ws.Range("N126").Value = -22482.25
$ws.Range("H131").Value = 3014.6726
$ws.Range("I131").Value = 458.66666
$ws.Range("J131").Value = 3973.175
$ws.Range("K131").Value = 1375.99998
$ws.Range("L131").Value = 11919.525
$ws.Range("M131").Value = 3664.00002
$ws.Range("N131").Value = -21999.525
$ws.Range("H132").Value = 2632.0789
$ws.Range("I132").Value = 2320.9524
$ws.Range("K132").Value = 20888.5716
$ws.Range("M132").Value = -18358.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 14859.44
$ws.Range("J136").Value = 12282
$ws.Range("L136").Value = 36846
$ws.Range("N136").Value = -41946
$ws.Range("H137").Value = 45620
$ws.Range("I137").Value = 20000
$ws.Range("J137").Value = 50744
$ws.Range("K137").Value = 20000
$ws.Range("L137").Value = 50744
$ws.Range("M137").Value = -14900
$ws.Range("N137").Value = -60944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4207.6875
$ws.Range("I7").Value = 3954.8667
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 3954.8667
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -3842.8667
$ws.Range("N7").Value = -8224
$ws.Range("H14").Value = 93336.664
$ws.Range("J14").Value = 93336.664
$ws.Range("L14").Value = 93336.664
$ws.Range("N14").Value = -93680.664
$ws.Range("H126").Value = 4207.6875
$ws.Range("I126").Value = 3954.8667
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 11864.6001
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -9394.6001
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 2815.889
$ws.Range("I132").Value = 2392.7727
$ws.Range("K132").Value = 7178.3181
$ws.Range("M132").Value = -4648.3181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 504709.6
$ws.Range("I14").Value = 3182.8333
$ws.Range("J14").Value = 1256999.8
$ws.Range("K14").Value = 3182.8333
$ws.Range("L14").Value = 1256999.8
$ws.Range("M14").Value = -3014.8333
$ws.Range("N14").Value = -1257335.8
$ws.Range("H100").Value = 624.7083
$ws.Range("I100").Value = 518.2308
$ws.Range("J100").Value = 750.5454999999999
$ws.Range("K100").Value = 1036.4616
$ws.Range("L100").Value = 1501.091
$ws.Range("M100").Value = -495.4616000000001
$ws.Range("N100").Value = -2583.091
$ws.Range("H126").Value = 1677.4
$ws.Range("I126").Value = 1771.75
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 5315.25
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -2845.25
$ws.Range("N126").Value = -8840
$ws.Range("H128").Value = 64355.832
$ws.Range("J128").Value = 64355.832
$ws.Range("L128").Value = 64355.832
$ws.Range("N128").Value = -74315.83199999999
$ws.Range("H136").Value = 3429.9783
$ws.Range("I136").Value = 2876.8438
$ws.Range("J136").Value = 4694.2856
$ws.Range("K136").Value = 8630.5314
$ws.Range("L136").Value = 14082.8568
$ws.Range("M136").Value = -6080.5314
$ws.Range("N136").Value = -19182.8568
$ws.Range("H139").Value = 63920.285
$ws.Range("J139").Value = 63920.285
$ws.Range("L139").Value = 63920.285
$ws.Range("N139").Value = -74200.285
